$d = $word.ActiveDocument

# 1) Update the letter date: September 19, 2025 -> September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the single-line mailing address into two paragraphs:
#    "121 9Th St, San Francisco CA 94103"
#      -> "121 9Th St"
#      -> "San Francisco, CA 94103"
$addrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "121 9Th St, San Francisco CA 94103*") {
        $addrPara = $p
        break
    }
}

if ($addrPara -ne $null) {
    $streetStart = $addrPara.Range.Start
    # "121 9Th St" is 10 characters long; split right after it.
    $streetRange = $d.Range($streetStart, $streetStart + 10)
    $streetRange.InsertParagraphAfter()

    # The remainder ", San Francisco CA 94103" now lives in its own paragraph;
    # rewrite it (formatting/rPr of the run is preserved) as "San Francisco, CA 94103".
    $cityPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like ", San Francisco CA 94103*") {
            $cityPara = $p
            break
        }
    }
    if ($cityPara -ne $null) {
        $cityPara.Range.Find.Execute(", San Francisco CA 94103", $false, $false, $false, $false, $false, $true, 1, $false, "San Francisco, CA 94103", 2) | Out-Null
    }
}

# 3) Remove the stray empty "No Spacing" paragraph directly below
#    "788 Minna Street Board of Directors".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $blank = $d.Paragraphs($i + 1)
        if ($blank.Range.Text.Trim() -eq "" -and $blank.Style.NameLocal -eq "No Spacing") {
            $blank.Range.Delete() | Out-Null
        }
        break
    }
}
